# Applies the recorded edit:
#  1. Rename header cell A1 from "URL" to "Google Maps Link".
#  2. Fill the blank "Phone" (column G) cells with a single space for the
#     rows that currently have no phone number recorded.
#  3. Turn on AutoFilter for the full data range (A1:H101), which also
#     registers the hidden "_xlnm._FilterDatabase" defined name.
#  4. Move the active selection from B1 to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "URL" header to "Google Maps Link".
$ws.Range("A1").Value2 = "Google Maps Link"

# 2. Rows whose Phone column (G) is currently empty get a single space.
$rowsMissingPhone = @(4, 9, 17, 30, 32, 36, 37, 38, 46, 48, 49, 50, 51, 52, 55, 61, 62, 64, 68, 69, 70, 71, 72, 75, 77, 82, 83, 84, 86, 89, 92, 100)
foreach ($r in $rowsMissingPhone) {
    $ws.Cells.Item($r, 7).Value2 = " "
}

# 3. Apply AutoFilter across the whole table and keep the filter-database
#    name hidden, the way Excel itself records it.
$dataRange = $ws.Range("A1:H101")
$dataRange.AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", $dataRange)
$filterName.Visible = $false

# 4. Move the selection to A1.
$ws.Range("A1").Select() | Out-Null
